$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-18 11:30:17"
$wsZhCn.Range("G5").Value = "2016-01-18 11:31:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-18 11:30:26"
$wsDeDe.Range("G5").Value = "2016-01-18 11:31:19"
